# Refresh the crypto Price (D) / Volume(1h) (E) columns, and the reshuffled
# coin rows at the bottom of the list (rows 48-51), per the GitHub Actions
# "cryptos list" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.491.97'
$ws.Range("E2").Value = '  -1.65%  '

$ws.Range("D3").Value = '2.591.02'
$ws.Range("E3").Value = '  -2.10%  '

$ws.Range("E4").Value = '  -0.34%  '

$ws.Range("D5").Value = '''561.75'
$ws.Range("E5").Value = '  -0.99%  '

$ws.Range("D6").Value = '''143.66'
$ws.Range("E6").Value = '  -2.42%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("E8").Value = '  -2.18%  '

$ws.Range("D9").Value = '2.600.25'
$ws.Range("E9").Value = '  -2.73%  '

$ws.Range("E10").Value = '  -2.54%  '

$ws.Range("E11").Value = '  -0.51%  '

$ws.Range("E12").Value = '  +10.15%  '

$ws.Range("E13").Value = '  +4.11%  '

$ws.Range("D14").Value = '3.049.96'
$ws.Range("E14").Value = '  -2.09%  '

$ws.Range("D15").Value = '''23.45'
$ws.Range("E15").Value = '  +7.03%  '

$ws.Range("D16").Value = '59.434.20'
$ws.Range("E16").Value = '  -1.76%  '

$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").Value = '2.598.45'
$ws.Range("E18").Value = '  -2.05%  '

$ws.Range("D19").Value = '''4.59'
$ws.Range("E19").Value = '  +0.52%  '

$ws.Range("D20").Value = '''339.11'
$ws.Range("E20").Value = '  -1.37%  '

$ws.Range("E21").Value = '  -0.39%  '

$ws.Range("D22").Value = '''6.54'
$ws.Range("E22").Value = '  +2.46%  '

$ws.Range("E23").Value = '  +0.34%  '

$ws.Range("D24").Value = '''63.77'
$ws.Range("E24").Value = '  -4.71%  '

$ws.Range("D25").Value = '''0.472'
$ws.Range("E25").Value = '  +7.18%  '

$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.48%  '

$ws.Range("D27").Value = '''0.161'
$ws.Range("E27").Value = '  -1.83%  '

$ws.Range("E28").Value = '  +1.09%  '

$ws.Range("D29").Value = '0.0₃0782'
$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("E31").Value = '  -2.06%  '

$ws.Range("E32").Value = '  -2.12%  '

$ws.Range("D33").Value = '''158.38'
$ws.Range("E33").Value = '  +2.13%  '

$ws.Range("E34").Value = '  -0.65%  '

$ws.Range("E35").Value = '  -0.72%  '

$ws.Range("E36").Value = '  +0.87%  '

$ws.Range("D37").Value = '''0.898'
$ws.Range("E37").Value = '  -1.31%  '

$ws.Range("E38").Value = '  -3.47%  '

$ws.Range("D39").Value = '''37.45'
$ws.Range("E39").Value = '  -0.32%  '

$ws.Range("E40").Value = '  -1.87%  '

$ws.Range("E41").Value = '  +0.67%  '

$ws.Range("D42").Value = '''295.01'
$ws.Range("E42").Value = '  -2.57%  '

$ws.Range("D43").Value = '''139.16'
$ws.Range("E43").Value = '  +8.40%  '

$ws.Range("D44").Value = '''0.998'
$ws.Range("E44").Value = '  +0.46%  '

$ws.Range("D45").Value = '''0.0978'
$ws.Range("E45").Value = '  -0.77%  '

$ws.Range("D46").Value = '''0.597'
$ws.Range("E46").Value = '  -1.61%  '

$ws.Range("D47").Value = '''10.63'
$ws.Range("E47").Value = '  -0.51%  '

$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '''0.0236'
$ws.Range("E48").Value = '  -0.14%  '

$ws.Range("B49").Value = 'Hedera'
$ws.Range("C49").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D49").Value = '''0.0532'
$ws.Range("E49").Value = '  -3.12%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''18.85'
$ws.Range("E50").Value = '  -0.71%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.963.45'
$ws.Range("E51").Value = '  -0.28%  '
